$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text: surgical character-range edits so the rest of the rich
# text (runs/fonts) for the cell is left alone.
# ---------------------------------------------------------------------
# A8 = "Volume 30   Number  17" -> "...  18" (only last digit changes)
$ws.Range("A8").Characters(22, 1).Text = "8"

# C9 = "Report Covering the Week  4/24/2023  Through  4/30/2023"
#   -> "...5/1/2023  Through  5/7/2023"
$ws.Range("C9").Characters(27, 9).Text = "5/1/2023"
$ws.Range("C9").Characters(46, 9).Text = "5/7/2023"

# ---------------------------------------------------------------------
# Cells whose TYPE flips between text ("0"/blank marker) and number.
# Copy the whole cell (value+style) from a same-shaped neighbor first
# so the style index lines up with the target (General vs #,##0),
# then stamp the real value on top (for the numeric ones).
# ---------------------------------------------------------------------
$ws.Range("C17").Copy($ws.Range("C15"))
$ws.Range("C15").Value = 1
$ws.Range("C17").Copy($ws.Range("F15"))
$ws.Range("F15").Value = 1
$ws.Range("C17").Copy($ws.Range("C16"))
$ws.Range("C16").Value = 5
$ws.Range("C17").Copy($ws.Range("C26"))
$ws.Range("C26").Value = 2
$ws.Range("C14").Copy($ws.Range("C27"))
$ws.Range("C14").Copy($ws.Range("C28"))
$ws.Range("C14").Copy($ws.Range("F28"))
$ws.Range("C14").Copy($ws.Range("C29"))
$ws.Range("C14").Copy($ws.Range("F29"))

# ---------------------------------------------------------------------
# Plain value updates (style/type unchanged).
# ---------------------------------------------------------------------
$ws.Range("M14").Value = 0
$ws.Range("H15").Value = -66.666666666666
$ws.Range("I15").Value = 3
$ws.Range("K15").Value = -40
$ws.Range("L15").Value = -62.5
$ws.Range("M15").Value = 50
$ws.Range("N15").Value = -70
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 25
$ws.Range("F16").Value = 14
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = 16.666666666666
$ws.Range("I16").Value = 50
$ws.Range("J16").Value = 67
$ws.Range("K16").Value = -25.373134328358
$ws.Range("L16").Value = 78.571428571428
$ws.Range("M16").Value = -44.444444444444
$ws.Range("N16").Value = -87.714987714987
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 150
$ws.Range("F17").Value = 33
$ws.Range("G17").Value = 25
$ws.Range("H17").Value = 32
$ws.Range("I17").Value = 134
$ws.Range("J17").Value = 109
$ws.Range("K17").Value = 22.935779816513
$ws.Range("L17").Value = 83.561643835616
$ws.Range("M17").Value = 88.732394366197
$ws.Range("N17").Value = 20.72072072072
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 300
$ws.Range("F18").Value = 16
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 77.777777777777
$ws.Range("I18").Value = 49
$ws.Range("J18").Value = 42
$ws.Range("K18").Value = 16.666666666666
$ws.Range("L18").Value = 22.5
$ws.Range("M18").Value = -44.943820224719
$ws.Range("N18").Value = -89.462365591397
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -33.333333333333
$ws.Range("G19").Value = 27
$ws.Range("H19").Value = -29.629629629629
$ws.Range("I19").Value = 103
$ws.Range("J19").Value = 135
$ws.Range("K19").Value = -23.703703703703
$ws.Range("L19").Value = 32.051282051282
$ws.Range("M19").Value = -8.849557522123
$ws.Range("N19").Value = -42.777777777777
$ws.Range("C20").Value = 3
$ws.Range("E20").Value = -25
$ws.Range("F20").Value = 17
$ws.Range("G20").Value = 18
$ws.Range("H20").Value = -5.555555555555
$ws.Range("I20").Value = 88
$ws.Range("J20").Value = 81
$ws.Range("K20").Value = 8.641975308641
$ws.Range("L20").Value = 54.38596491228
$ws.Range("M20").Value = 6.024096385542
$ws.Range("N20").Value = -92.93172690763
$ws.Range("C21").Value = 29
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = 31.818181818181
$ws.Range("F21").Value = 100
$ws.Range("G21").Value = 94
$ws.Range("H21").Value = 6.382978723404
$ws.Range("I21").Value = 429
$ws.Range("J21").Value = 440
$ws.Range("K21").Value = -2.5
$ws.Range("L21").Value = 49.477351916376
$ws.Range("M21").Value = -4.666666666666
$ws.Range("N21").Value = -82.287365813377
$ws.Range("L22").Value = 100
$ws.Range("C24").Value = 30
$ws.Range("D24").Value = 33
$ws.Range("E24").Value = -9.090909090909
$ws.Range("F24").Value = 117
$ws.Range("G24").Value = 120
$ws.Range("H24").Value = -2.5
$ws.Range("I24").Value = 519
$ws.Range("J24").Value = 481
$ws.Range("K24").Value = 7.900207900207
$ws.Range("L24").Value = 45.378151260504
$ws.Range("M24").Value = 94.38202247191
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 83.333333333333
$ws.Range("F25").Value = 39
$ws.Range("G25").Value = 33
$ws.Range("H25").Value = 18.181818181818
$ws.Range("I25").Value = 183
$ws.Range("J25").Value = 164
$ws.Range("K25").Value = 11.585365853658
$ws.Range("L25").Value = 42.96875
$ws.Range("M25").Value = -8.040201005025
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = -25
$ws.Range("I26").Value = 9
$ws.Range("K26").Value = -25
$ws.Range("L26").Value = -25
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -100
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -25
$ws.Range("J27").Value = 17
$ws.Range("K27").Value = 11.764705882352
$ws.Range("H28").Value = -100
$ws.Range("I28").Value = 4
$ws.Range("K28").Value = -42.857142857142
$ws.Range("L28").Value = -20
$ws.Range("M28").Value = 33.333333333333
$ws.Range("N28").Value = -71.428571428571
$ws.Range("H29").Value = -100
$ws.Range("I29").Value = 4
$ws.Range("K29").Value = -20
$ws.Range("L29").Value = 33.333333333333
$ws.Range("M29").Value = 33.333333333333
$ws.Range("N29").Value = -60
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 0
